# Atualiza base de usuários - adiciona novo usuário "carlos" com senha "2020caca"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("A5").Value = "carlos"
$ws.Range("B5").Value = "2020caca"

$ws.Range("B6").Select()
